$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 11;  D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 15;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 31;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 32;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 36;  D = "[1, 1, 0, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment']" },
    @{ Row = 43;  D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 45;  D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 54;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 58;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 61;  D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 68;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 81;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 84;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 92;  D = "[1, 0, 1, 0, 0, 0, 1]"; E = "['Normal', 'HardwareFault', 'SoftwareFault']" },
    @{ Row = 97;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 109; D = "[1, 1, 0, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment']" },
    @{ Row = 116; D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" }
)

foreach ($u in $updates) {
    $ws.Range("D$($u.Row)").Value = $u.D
    $ws.Range("E$($u.Row)").Value = $u.E
}
